{"js": "const body = context.document.body;\n\n// ---------------------------------------------------------------------\n// 1. \"In terms of flexibility...\" paragraph:\n//    - \"an interface for drinks\" -> \"an abstract class for drinks\"\n//    - \"...These new additions would simply be realizations that implement\n//       the drink interface.\" ->\n//      \"...These new additions would simply be added through inheritance\n//       that are subclasses of the drink abstract class.\"\n//    - the empty paragraph that used to follow is removed (merged away)\n// ---------------------------------------------------------------------\nlet results = body.search(\"an interface for drinks\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nresults.items[0].insertText(\"an abstract class for drinks\", Word.InsertLocation.replace);\nawait context.sync();\n\nresults = body.search(\n  \"as well as ice or artificial sweeteners. These new additions would simply be realizations that implement the drink interface.\",\n  { matchCase: true }\n);\nresults.load(\"items\");\nawait context.sync();\nresults.items[0].insertText(\n  \"as well as ice or artificial sweeteners. These new additions would simply be added through inheritance that are subclasses of the drink abstract class.\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 2. Remove the empty paragraph following the \"flexibility\" paragraph.\n// ---------------------------------------------------------------------\nlet paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfunction findEmptyParaAfter(items, searchStartIdx) {\n  for (let i = searchStartIdx; i < items.length; i++) {\n    if (items[i].text === \"\") {\n      return i;\n    }\n  }\n  return -1;\n}\n\nlet idx = findEmptyParaAfter(paragraphs.items, 0);\nparagraphs.items[idx].delete();\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 2b. \"Simplicity and \" + \"Understandability\" runs merge into a single run.\n// ---------------------------------------------------------------------\nresults = body.search(\"Simplicity and Understandability\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nresults.items[0].insertText(\"Simplicity and Understandability\", Word.InsertLocation.replace);\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 3. \"that it was not cluttered, easy for developers to read, and simple\n//     enough to understand.\" paragraph gets extended, and the empty\n//     paragraph that used to follow is removed.\n// ---------------------------------------------------------------------\nresults = body.search(\n  \"that it was not cluttered, easy for developers to read, and simple enough to understand.\",\n  { matchCase: true }\n);\nresults.load(\"items\");\nawait context.sync();\nresults.items[0].insertText(\n  \"that it was not cluttered, easy for developers to read, and simple enough to understand. The abstract class and all the inheritance should separate everything in such a way that it will be easy to navigate and understand everything\\u2019s purpose.\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\nparagraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\nidx = findEmptyParaAfter(paragraphs.items, 0);\nparagraphs.items[idx].delete();\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 4. \"I wanted to use a condiments and drink interface...\" paragraph:\n//    \"individual realizations of their classes.\" -> \"individual\n//    inheritances of their classes.\"; empty paragraph after it removed.\n// ---------------------------------------------------------------------\nresults = body.search(\"individual realizations of their classes.\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nresults.items[0].insertText(\"individual inheritances of their classes.\", Word.InsertLocation.replace);\nawait context.sync();\n\nparagraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\nidx = findEmptyParaAfter(paragraphs.items, 0);\nparagraphs.items[idx].delete();\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 5. \"I utilized realization while using a drink interface, ...\" ->\n//    \"I utilized inheritance while using a drink abstract class, ...\"\n// ---------------------------------------------------------------------\nresults = body.search(\n  \"I utilized realization while using a drink interface, and the machine itself will be \\u201ccomposed\\u201d of all the possible options\",\n  { matchCase: true }\n);\nresults.load(\"items\");\nawait context.sync();\nif (results.items.length === 0) {\n  // Fall back in case of straight quotes in the document instead of curly ones.\n  results = body.search(\n    \"I utilized realization while using a drink interface, and the machine itself will be \\\"composed\\\" of all the possible options\",\n    { matchCase: true }\n  );\n  results.load(\"items\");\n  await context.sync();\n}\nresults.items[0].insertText(\n  \"I utilized inheritance while using a drink abstract class, and the machine itself will be \\\"composed\\\" of all the possible options\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 6. Remove the empty paragraph following \"so that it was easily readable\n//    and understandable.\"\n// ---------------------------------------------------------------------\nparagraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\nidx = findEmptyParaAfter(paragraphs.items, 0);\nparagraphs.items[idx].delete();\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 7. \"I will have several packages: ...\" paragraph gets simplified\n//    (proofing-error markers removed, \"interfaces and realizations\" ->\n//    \"abstract classes and implementations\", \"max\" + \"quantities\" merged).\n//    We replace the *entire* paragraph range (not just a sub-range) so\n//    that the stray w:proofErr markers (gramStart/gramEnd/spellStart/\n//    spellEnd) scattered between the original runs are swept away too.\n// ---------------------------------------------------------------------\nparagraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet packagesPara = null;\nfor (const p of paragraphs.items) {\n  if (p.text.indexOf(\"I will have several packages\") === 0) {\n    packagesPara = p;\n    break;\n  }\n}\n\nconst packagesParaRange = packagesPara.getRange();\npackagesParaRange.insertText(\n  \"I will have several packages: beverage(for the drink and condiment abstract class and all implementations) and machine (for the DrinkMachine class). The abstract classes and implementations will be implemented as described above with associated prices, max quantities for condiments, etc. The machine will run the entire program such as asking the user for their drink preference, drink type, and condiment choices then dispense the drink. All of these will have associated JUnit tests.\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Get-FirstEmptyParagraphIndex($doc) {\n    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {\n        $p = $doc.Paragraphs.Item($i)\n        $t = $p.Range.Text -replace \"`r\", \"\"\n        if ($t -eq \"\") {\n            return $i\n        }\n    }\n    return -1\n}\n\n# Locates $searchText anywhere in the document and replaces it with\n# $replaceText. When $replaceText contains a straight double-quote\n# character, the found range's .Text is overwritten directly (rather than\n# going through Find's Replacement argument) so that Word's smart-quote\n# autocorrect does not turn it into a curly quote; otherwise Find's own\n# Replacement mechanism is used because it also merges same-formatted\n# adjacent runs the way a normal Word edit would.\nfunction Replace-Text($doc, $searchText, $replaceText) {\n    $range = $doc.Content\n    $find = $range.Find\n    if ($replaceText.Contains('\"')) {\n        $found = $find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false)\n        $range.Text = $replaceText\n    } else {\n        $find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n    }\n}\n\n# Same idea, but scoped to a single paragraph's range - used to replace an\n# entire paragraph (sweeping away any stray w:proofErr markers left behind\n# between the original runs).\nfunction Replace-ParagraphText($paragraph, $replaceText) {\n    $range = $paragraph.Range\n    $find = $range.Find\n    if ($replaceText.Contains('\"')) {\n        $found = $find.Execute($range.Text, $false, $false, $false, $false, $false, $true, 1, $false)\n        $range.Text = $replaceText\n    } else {\n        $find.Execute($range.Text, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n    }\n}\n\n# ---------------------------------------------------------------------\n# 1. \"In terms of flexibility...\" paragraph:\n#    - \"an interface for drinks\" -> \"an abstract class for drinks\"\n#    - \"...These new additions would simply be realizations that implement\n#       the drink interface.\" ->\n#      \"...These new additions would simply be added through inheritance\n#       that are subclasses of the drink abstract class.\"\n#    - the empty paragraph that used to follow is removed (merged away)\n# ---------------------------------------------------------------------\nReplace-Text $d \"an interface for drinks\" \"an abstract class for drinks\"\n\nReplace-Text $d `\n    \"as well as ice or artificial sweeteners. These new additions would simply be realizations that implement the drink interface.\" `\n    \"as well as ice or artificial sweeteners. These new additions would simply be added through inheritance that are subclasses of the drink abstract class.\"\n\n$idx = Get-FirstEmptyParagraphIndex $d\n$d.Paragraphs.Item($idx).Range.Delete()\n\n# ---------------------------------------------------------------------\n# 2. \"Simplicity and \" + \"Understandability\" runs merge into a single run.\n# ---------------------------------------------------------------------\nReplace-Text $d \"Simplicity and Understandability\" \"Simplicity and Understandability\"\n\n# ---------------------------------------------------------------------\n# 3. \"that it was not cluttered, easy for developers to read, and simple\n#     enough to understand.\" paragraph gets extended, and the empty\n#     paragraph that used to follow is removed.\n# ---------------------------------------------------------------------\nReplace-Text $d `\n    \"that it was not cluttered, easy for developers to read, and simple enough to understand.\" `\n    \"that it was not cluttered, easy for developers to read, and simple enough to understand. The abstract class and all the inheritance should separate everything in such a way that it will be easy to navigate and understand everything\u2019s purpose.\"\n\n$idx = Get-FirstEmptyParagraphIndex $d\n$d.Paragraphs.Item($idx).Range.Delete()\n\n# ---------------------------------------------------------------------\n# 4. \"I wanted to use a condiments and drink interface...\" paragraph:\n#    \"individual realizations of their classes.\" -> \"individual\n#    inheritances of their classes.\"; empty paragraph after it removed.\n# ---------------------------------------------------------------------\nReplace-Text $d \"individual realizations of their classes.\" \"individual inheritances of their classes.\"\n\n$idx = Get-FirstEmptyParagraphIndex $d\n$d.Paragraphs.Item($idx).Range.Delete()\n\n# ---------------------------------------------------------------------\n# 5. \"I utilized realization while using a drink interface, ...\" ->\n#    \"I utilized inheritance while using a drink abstract class, ...\"\n# ---------------------------------------------------------------------\nReplace-Text $d `\n    'I utilized realization while using a drink interface, and the machine itself will be \"composed\" of all the possible options' `\n    'I utilized inheritance while using a drink abstract class, and the machine itself will be \"composed\" of all the possible options'\n\n# ---------------------------------------------------------------------\n# 6. Remove the empty paragraph following \"so that it was easily readable\n#    and understandable.\"\n# ---------------------------------------------------------------------\n$idx = Get-FirstEmptyParagraphIndex $d\n$d.Paragraphs.Item($idx).Range.Delete()\n\n# ---------------------------------------------------------------------\n# 7. \"I will have several packages: ...\" paragraph gets simplified\n#    (proofing-error markers removed, \"interfaces and realizations\" ->\n#    \"abstract classes and implementations\", \"max\" + \"quantities\" merged).\n#    We replace the *entire* paragraph range (not just a sub-range) so\n#    that the stray proofing-error markers scattered between the original\n#    runs are swept away too.\n# ---------------------------------------------------------------------\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text.StartsWith(\"I will have several packages\")) {\n        Replace-ParagraphText $p \"I will have several packages: beverage(for the drink and condiment abstract class and all implementations) and machine (for the DrinkMachine class). The abstract classes and implementations will be implemented as described above with associated prices, max quantities for condiments, etc. The machine will run the entire program such as asking the user for their drink preference, drink type, and condiment choices then dispense the drink. All of these will have associated JUnit tests.\"\n        break\n    }\n}\n"}
